$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date for first row.
# This timestamp is shared with de-de's Correspond Handoff Datetime for the
# same source file (same underlying handoff event), so both must be updated
# together to keep them in sync.
$overview.Range("G2").Value = "2016-08-19 03:01:58"
$dede.Range("H2").Value = "2016-08-19 03:01:58"

# zh-cn sheet row 2: Correspond Handoff Datetime / Correspond Handback DateTime
$zhcn.Range("H2").Value = "2016-08-19 03:01:52"
$zhcn.Range("K2").Value = "2016-08-19 03:02:15"

# de-de sheet row 2: Correspond Handback DateTime
$dede.Range("K2").Value = "2016-08-19 03:02:22"
